$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.935.94'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.389.45'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.32'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.08'
$ws.Range("E6").Value = '  +1.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.380.79'
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.169'
$ws.Range("E10").Value = '  +11.78%  '
$ws.Range("E11").Value = '  +3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.55'
$ws.Range("E12").Value = '  +2.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000279'
$ws.Range("E13").Value = '  +5.57%  '
$ws.Range("E14").Value = '  +3.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.929.29'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.29'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.389.02'
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '64.834.36'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.82'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.993'
$ws.Range("E21").Value = '  +2.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '471.37'
$ws.Range("E22").Value = '  +14.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.98'
$ws.Range("E23").Value = '  +13.44%  '
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.09'
$ws.Range("E25").Value = '  +5.14%  '
$ws.Range("E26").Value = '  -1.70%  '
$ws.Range("E27").Value = '  +7.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.80'
$ws.Range("E28").Value = '  +2.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.77'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.00'
$ws.Range("E30").Value = '  +7.05%  '
$ws.Range("E31").Value = '  +4.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.52'
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '572.11'
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.44'
$ws.Range("E34").Value = '  +6.37%  '
$ws.Range("E35").Value = '  +2.16%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.55'
$ws.Range("E37").Value = '  +4.61%  '
$ws.Range("E38").Value = '  -4.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.65'
$ws.Range("E39").Value = '  +2.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0750'
$ws.Range("E40").Value = '  +2.28%  '
$ws.Range("E41").Value = '  +1.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.091.64'
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  +2.92%  '
$ws.Range("E45").Value = '  +4.31%  '
$ws.Range("E46").Value = '  +5.91%  '
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '138.09'
$ws.Range("E50").Value = '  +4.00%  '
$ws.Range("E51").Value = '  +4.25%  '
